# DetailedLog.xlsx update: log updated, KNN R file added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------
# 1. Row 19: add the "Start KNN research" note in column G
# ---------------------------------------------------------------
$ws.Range("G19").Value = "Start KNN research"

# ---------------------------------------------------------------
# 2. New data rows 20-24 (KNN research log entries)
# ---------------------------------------------------------------

# -- Row 20 --
$ws.Range("A20").Value = 45706
$ws.Range("B20").Value = "Started research for the KNN model what what I want to achieve with it. Did a simple model with the famous iris dataset. Wanted to do it on a harder dataset so chose a wine dataset to predict the quality of the wine"
$ws.Range("C20").Value = "Done"
$ws.Range("D20").Value = "Medium"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Wine dataset downloaded in csv. Business objectives started"
$ws.Range("G20").Value = "Train model KNN"

# -- Row 21 --
$ws.Range("A21").Value = 45706
$ws.Range("B21").Value = "Simple model trained"
$ws.Range("C21").Value = "Done"
$ws.Range("D21").Value = "Medium"
$ws.Range("E21").Value = 0.3
$ws.Range("F21").Value = "Poor accuracy of 51%"
$ws.Range("G21").Value = "complete data visualisation"

# -- Row 22 --
$ws.Range("A22").Value = 45706
$ws.Range("B22").Value = "Data visualisation to give insights into the data"
$ws.Range("C22").Value = "Done"
$ws.Range("D22").Value = "Medium"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = "Produced 3 data visualisation graphs. Distribution, acidity by quality , content by quality"
$ws.Range("G22").Value = "Improve model"

# -- Row 23 --
$ws.Range("A23").Value = 45707
$ws.Range("B23").Value = "Improve model by cleaning data, data cleaning, removing columns with no correlation, finding the best k value"
$ws.Range("C23").Value = "Done"
$ws.Range("D23").Value = "Medium"
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = "Improved to 59% through tunning, still poor accuracy"
$ws.Range("G23").Value = "Documentation on portfolio"

# -- Row 24 --
$ws.Range("A24").Value = 45707
$ws.Range("B24").Value = "Documentation on porfolio of process"
$ws.Range("C24").Value = "Done"
$ws.Range("D24").Value = "Medium"
$ws.Range("E24").Value = 3

# ---------------------------------------------------------------
# 3. Formatting for rows 20-24
#    column A -> same date/center/wrap style already used by A12:A19
#    columns B:G -> same center/wrap style already used across the table
# ---------------------------------------------------------------
$ws.Range("A19").Copy()
$ws.Range("A20:A24").PasteSpecial($xlPasteFormats)

$ws.Range("B19").Copy()
$ws.Range("B20:G24").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false

# Row heights (auto-fit heights for the wrapped multi-line rows)
$ws.Rows.Item(20).RowHeight = 57.6
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 28.8

# ---------------------------------------------------------------
# 4. Extra trailing blank rows (25-31)
#    25-28 keep the regular center/wrap style
#    29-31 use a new center-only (no wrap) style
# ---------------------------------------------------------------
$ws.Range("B19").Copy()
$ws.Range("A25:G28").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("A29:G31").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------
# 5. Update the active selection to F22
# ---------------------------------------------------------------
$ws.Range("F22").Select() | Out-Null
